$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.840.02"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.54%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.873.34"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.86%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "301.46"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.11%  "
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5343"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.61%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3746"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07198"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.67%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.61"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.33%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8900"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.84%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08184"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.22%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.881.52"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +6.51%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "93.13"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.45%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.307"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.28%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.003"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.14%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.83"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.36%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008535"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.84%  "
$ws.Range("E19").Value = "  +0.03%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "26.873.25"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.53%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.990"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.80%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.62"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.96%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.383"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.91%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.284"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.56%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.07"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.68%  "
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.07"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.05%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "114.02"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.48%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.709"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.15%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.617"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.43%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09114"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.61%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.8110"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.80%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05015"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.39%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.174"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.47%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.961"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.63%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6108"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.87%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.657"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.14%  "
$ws.Range("E38").Value = "  -4.56%  "
$ws.Range("E39").Value = "  -2.42%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.068"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.47%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.579"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.47%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.883"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.25%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5173"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.90%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "114.91"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.85%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1496"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.88%  "
$ws.Range("E46").Value = "  +0.01%  "
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.644"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.07%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.944"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.44%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "37.54"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.04%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06057"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.04%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "62.27"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.67%  "
